# CORP_holdings.xlsx update:
#   - Bump the "as of" date in the confidential disclaimer banner
#     (2021-05-14 -> 2021-05-17)
#   - Refresh the Weight (D) / Percent Change (E) figures for the
#     current holdings table (rows 2-9)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet ships protected; unprotect so the cells below can be written,
# then restore protection at the end.
$ws.Unprotect()

# --- Disclaimer banner text (A12) -----------------------------------
$ws.Range("A12").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-17 for illustrative purposes only and are subject to change."

# --- Holdings table: Weight / Percent Change refresh ------------------
$ws.Range("D2").Value = 0.1774293542387457
$ws.Range("E2").Value = 0

$ws.Range("D3").Value = 0.1772953622782633
$ws.Range("E3").Value = 0

$ws.Range("D4").Value = 0.2255084694918305
$ws.Range("E4").Value = 0

$ws.Range("D5").Value = 0.07986120832750035
$ws.Range("E5").Value = 0

$ws.Range("D6").Value = 0.079667219966802
$ws.Range("E6").Value = 0

$ws.Range("D7").Value = 0.1202347859128452
$ws.Range("E7").Value = -0.0009813542688908994

$ws.Range("D8").Value = 0.140003599784013
$ws.Range("E8").Value = 0

$ws.Range("E9").Value = -0.0001179929204245811

# Restore the original sheet protection state.
$ws.Protect()
